$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the updated numeric-looking cells to remain plain text (matching the
# workbook's existing inlineStr/text convention for these columns) instead of
# being auto-converted to numbers/percentages by Excel's input parser.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.22%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.56%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.209"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.80%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07670"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.60%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.632"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.62%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9161"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.35%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.25%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "9.67%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1825"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.36%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09108"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.85%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.54%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001257"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.25%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005737"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.71%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.340"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.41%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.302"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.26%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3335"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.31%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.377"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "11.96%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1381"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.24%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.72%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.08%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004383"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.88%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.07%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02496"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.07%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05342"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.06%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007850"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.08%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1314"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.89%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006502"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.21%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001914"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.93%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008255"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.20%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3330"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.05%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006706"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.37%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.12%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2702"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "768.44%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-26.23%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.12%"
